$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Resultados")

$ws.Range("I4").Value = 152600
$ws.Range("E5").Value = 8
$ws.Range("I5").Value = 154800
$ws.Range("E7").Value = 9
$ws.Range("I7").Value = 210050
$ws.Range("G11").Value = 10.17333333333333
$ws.Range("G12").Value = 5.16
$ws.Range("G14").Value = 5.25125
